$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: remove the "_GoBack" bookmark that currently sits right
# after "March 29" (between "March 29" and ", 2017").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# Change 2: split the "[  ] Edit questions and answers from table"
# run (in the "From Last Week" checklist) into separate runs with
# proofErr (grammar-check) markers bracketing "[  ]".
#   <w:r><w:tab/></w:r>
#   <w:proofErr w:type="gramStart"/>
#   <w:r><w:t>[  ]</w:t></w:r>
#   <w:proofErr w:type="gramEnd"/>
#   <w:r><w:t xml:space="preserve"> Edit questions and answers from table</w:t></w:r>
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("[  ] Edit questions and answers from table")) {
        $targetPara = $p
    }
}

if ($targetPara -ne $null) {
    $r = $targetPara.Range
    $xmlFrag = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="05316CD6" w14:textId="4409940D" w:rsidR="00881070" w:rsidRDefault="00881070" w:rsidP="00881070"><w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>[  ]</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Edit questions and answers from table</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
    $r.InsertXML($xmlFrag)
}

# ------------------------------------------------------------------
# Change 3: append the commit URL right after "GitHub Link: " and
# re-create the "_GoBack" bookmark so it ends up immediately after
# the newly-inserted link text.
# ------------------------------------------------------------------
$githubPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("GitHub Link:")) {
        $githubPara = $p
    }
}

if ($githubPara -ne $null) {
    $pr = $githubPara.Range
    # Insert right before the paragraph mark (pr.End is just past the
    # paragraph's final character, so pr.End - 1 is the paragraph-mark
    # position itself; inserting there places the text before it).
    $insertAt = $d.Range($pr.End - 1, $pr.End - 1)
    $url = "https://github.com/aroja1/Program-Project-Seminar-for-Minors/commit/b092c4811043c4a4b48cad3382fccfd2acc80a42"
    $insertAt.InsertAfter($url)

    # Re-read the paragraph's end now that the link has been inserted.
    $afterLink = $githubPara.Range.End - 1

    # A collapsed range sitting exactly on the last character position
    # of a paragraph confuses Bookmarks.Add in this engine, so nudge
    # it out of that special spot with a temporary placeholder
    # character, add the bookmark, then remove the placeholder again.
    $placeholderRange = $d.Range($afterLink, $afterLink)
    $placeholderRange.InsertAfter("X")

    $bookmarkRange = $d.Range($afterLink, $afterLink)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)

    $placeholderDelRange = $d.Range($afterLink, $afterLink + 1)
    $placeholderDelRange.Delete()
}
